$wb = $excel.ActiveWorkbook

$oldVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)"

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = "Version: " + $newVersion
$wsAbout.Range("A6").Value = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Sanyuan Coal Mine, China, M2103, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($r = 2; $r -le 8; $r++) {
    $cell = $wsData.Cells.Item($r, 19)
    if ($cell.Value() -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
